# "updated legacy GSC export data"
# The oldest day (2025-11-03) was dropped from the export, so its row is
# removed from the data table and every subsequent row shifts up one
# position. The two most recent days in the refreshed export
# (2025-11-04 and 2025-11-05, now sitting in rows 2 and 3) don't have
# "No video indexed" / "Video indexed" counts yet, so those cells are
# left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for 2025-11-03; everything below shifts up one row.
$ws.Rows("2:2").Delete()

# The newest two dates don't have indexing counts yet - blank them out.
$ws.Range("B2:C3").ClearContents()
